# Applies the commit's data changes:
#  - RUNMANAGER: row3 "execute" Yes->yes; new row4 "invalidDeliveryAddress" test case
#  - DATA: rows 2-4 tweaks (browser chrome->firefox on two rows, items/first/last/zip
#    filled in on row4), new row5 (addItemsToCartTest w/ chrome + 2 items) and new
#    row6 (invalidDeliveryAddress data row)
#  - selection moved in both sheets (C3 on RUNMANAGER, A9 on DATA) with DATA staying
#    the active/tab-selected sheet

$wb = $excel.ActiveWorkbook

$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$data = $wb.Worksheets.Item("DATA")

# ---- RUNMANAGER sheet ----
$runmanager.Range("C3").Value = "yes"

$runmanager.Range("A4").Value = "invalidDeliveryAddress"
$runmanager.Range("B4").Value = "Test invalid delivery details"
$runmanager.Range("C4").Value = "no"
$runmanager.Range("D4").Value = "1"
$runmanager.Range("E4").Value = "1"

# ---- DATA sheet ----
$data.Range("C3").Value = "firefox"

$data.Range("C4").Value = "firefox"
$data.Range("F4").Value = "Sauce Labs Fleece Jacket"
$data.Range("G4").Value = "John"
$data.Range("H4").Value = "Wick"
# I4/I5 are digit strings ("1010"/"0628") that must stay text (not be coerced
# to numbers, which would also eat the leading zero in "0628"): force a text
# format before assigning, then drop back to the workbook's normal style.
$data.Range("I4").NumberFormat = "@"
$data.Range("I4").Value = "1010"
$data.Range("I4").Style = "Normal"

$data.Range("A5").Value = "addItemsToCartTest"
$data.Range("B5").Value = "yes"
$data.Range("C5").Value = "chrome"
$data.Range("D5").Value = "standard_user"
$data.Range("E5").Value = "secret_sauce"
$data.Range("F5").Value = "Sauce Labs Fleece Jacket;Test.allTheThings() T-Shirt (Red)"
$data.Range("G5").Value = "Ethan"
$data.Range("H5").Value = "Hunt"
$data.Range("I5").NumberFormat = "@"
$data.Range("I5").Value = "0628"
$data.Range("I5").Style = "Normal"

$data.Range("A6").Value = "invalidDeliveryAddress"
$data.Range("B6").Value = "no"
$data.Range("C6").Value = "chrome"
$data.Range("D6").Value = "standard_user"
$data.Range("E6").Value = "secret_sauce"

# The G4:I4 / G5:I5 / F6:I6 cells above are blank "item/firstname/lastname/zip
# not applicable" placeholders that carry a quote-prefixed blank style in the
# original sheet (see F2:I3). Writing through `.Value` on a previously-blank
# styled cell (G4:I5) -- or creating a brand new cell (row 6) -- does not
# carry that style along, so restore it explicitly via a format-only paste
# from an existing cell that already has it.
$data.Range("I2").Copy()
$data.Range("G4:I4").PasteSpecial(-4122)
$data.Range("G5:I5").PasteSpecial(-4122)
$data.Range("F6:I6").PasteSpecial(-4122)

# ---- selection / active sheet bookkeeping ----
$runmanager.Activate()
$runmanager.Range("C3").Select()

$data.Activate()
$data.Range("A9").Select()
